$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 7 that duplicates row 6 (same values and formatting:
# 26-09-2025 date in A and its corresponding gold price text in B).
$ws.Range("A6:B6").Copy()
$ws.Range("A7").PasteSpecial()
